$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '66.228.00'
$ws.Range("E2").Value = '  -0.84%  '

# Row 3
$ws.Range("D3").Value = '3.318.23'
$ws.Range("E3").Value = '  -1.42%  '

# Row 4
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '188.64'
$ws.Range("E5").Value = '  +2.69%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '559.78'
$ws.Range("E6").Value = '  -0.24%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.13%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.586'
$ws.Range("E8").Value = '  -1.52%  '

# Row 9
$ws.Range("D9").Value = '3.310.59'
$ws.Range("E9").Value = '  -1.44%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.186'
$ws.Range("E10").Value = '  -0.97%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.587'
$ws.Range("E11").Value = '  -1.52%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.70'
$ws.Range("E12").Value = '  -0.35%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000272'
$ws.Range("E13").Value = '  +1.53%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.64'
$ws.Range("E14").Value = '  -1.16%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '633.27'
$ws.Range("E15").Value = '  +3.47%  '

# Row 16
$ws.Range("D16").Value = '3.845.30'
$ws.Range("E16").Value = '  -1.52%  '

# Row 17
$ws.Range("E17").Value = '  +2.88%  '

# Row 18
$ws.Range("D18").Value = '66.180.75'
$ws.Range("E18").Value = '  -0.73%  '

# Row 19
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.343.91'
$ws.Range("E19").Value = '  -0.87%  '

# Row 20
$ws.Range("B20").Value = 'TRON'
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.117'
$ws.Range("E20").Value = '  -1.21%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.10'
$ws.Range("E21").Value = '  -4.98%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.910'
$ws.Range("E22").Value = '  -0.32%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.08'
$ws.Range("E23").Value = '  +6.35%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '103.38'
$ws.Range("E24").Value = '  +7.79%  '

# Row 25
$ws.Range("E25").Value = '  -2.33%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.95'
$ws.Range("E26").Value = '  -3.60%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.04'
$ws.Range("E27").Value = '  +0.72%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.73'
$ws.Range("E28").Value = '  -0.73%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.59'
$ws.Range("E29").Value = '  +0.22%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.69'
$ws.Range("E30").Value = '  -1.25%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.34'
$ws.Range("E31").Value = '  -1.53%  '

# Row 32
$ws.Range("E32").Value = '  +5.86%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.38'
$ws.Range("E33").Value = '  +0.36%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.09'
$ws.Range("E34").Value = '  -1.50%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '550.57'
$ws.Range("E35").Value = '  +3.44%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.105'
$ws.Range("E36").Value = '  -0.03%  '

# Row 37
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '3.838.38'
$ws.Range("E37").Value = '  +1.97%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '57.73'
$ws.Range("E38").Value = '  -1.55%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.20%  '

# Row 40
$ws.Range("D40").Value = '0.0₃0736'
$ws.Range("E40").Value = '  +1.59%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.31'
$ws.Range("E41").Value = '  -3.24%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '33.78'
$ws.Range("E42").Value = '  +3.27%  '

# Row 43
$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.70'
$ws.Range("E43").Value = '  -0.62%  '

# Row 44
$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.128'
$ws.Range("E44").Value = '  +0.33%  '

# Row 45
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.336'
$ws.Range("E45").Value = '  -4.69%  '

# Row 46
$ws.Range("B46").Value = 'CoreDAO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.21'
$ws.Range("E46").Value = '  -14.40%  '

# Row 47
$ws.Range("E47").Value = '  +0.50%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.26'
$ws.Range("E48").Value = '  +2.76%  '

# Row 49
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.129'
$ws.Range("E49").Value = '  -0.71%  '

# Row 50
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.61'
$ws.Range("E50").Value = '  -2.99%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.998'
